$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 - this shifts existing rows 37..117 down to 38..118
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with the new record
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 45028
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = "Tropicales y subtropicales"
$ws.Range("I37").Value = 100108002
$ws.Range("J37").Value = "Mango"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 60
$ws.Range("N37").Value = 8000
$ws.Range("O37").Value = 8500
$ws.Range("P37").Value = 8250
$ws.Range("Q37").Value = '$/bandeja 4 kilos'
$ws.Range("R37").Value = "Perú"
$ws.Range("S37").Value = 2062
$ws.Range("T37").Value = 4
